$wb = $excel.ActiveWorkbook

# --- Update AccountCreationData sheet (email addresses for rows 2-4) ---
$ws = $wb.Worksheets.Item("AccountCreationData")
$ws.Range("A2").Value = "newtest4@gmail.com"
$ws.Range("A3").Value = "newtest5@gmail.com"
$ws.Range("A4").Value = "newtest6@gmail.com"

# Move the selection/active cell on this sheet to B13
$ws.Range("B13").Select()

# Make AccountCreationData the active (selected) sheet/tab
$ws.Activate()
